$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.180.90"
$ws.Range("E2").Value = "  -2.64%  "
$ws.Range("D3").Value = "1.821.18"
$ws.Range("E3").Value = "  -2.12%  "
$ws.Range("D4").Value = "'1.006"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -1.20%  "
$ws.Range("D5").Value = "'314.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.23%  "
$ws.Range("D6").Value = "'1.006"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.13%  "
$ws.Range("D7").Value = "'0.4249"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.32%  "
$ws.Range("D8").Value = "'0.3677"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.15%  "
$ws.Range("D9").Value = "'0.07240"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.73%  "
$ws.Range("D10").Value = "'0.8605"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.11%  "
$ws.Range("E11").Value = "  -3.81%  "
$ws.Range("D12").Value = "1.823.00"
$ws.Range("E12").Value = "  -2.01%  "
$ws.Range("D13").Value = "'6.700"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.56%  "
$ws.Range("D14").Value = "'0.07092"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.95%  "
$ws.Range("D15").Value = "'5.303"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.72%  "
$ws.Range("D16").Value = "'87.93"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.45%  "
$ws.Range("E17").Value = "  -1.36%  "
$ws.Range("D18").Value = "'0.000008862"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.18%  "
$ws.Range("D19").Value = "'1.005"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.09%  "
$ws.Range("D20").Value = "'15.03"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.35%  "
$ws.Range("D21").Value = "27.217.84"
$ws.Range("E21").Value = "  -2.52%  "
$ws.Range("D22").Value = "'5.136"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.79%  "
$ws.Range("E23").Value = "  -3.43%  "
$ws.Range("D24").Value = "2.046.16"
$ws.Range("E24").Value = "  -1.48%  "
$ws.Range("D25").Value = "'2.001"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.44%  "
$ws.Range("D26").Value = "'153.35"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.51%  "
$ws.Range("D27").Value = "'18.28"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.46%  "
$ws.Range("D28").Value = "'2.127"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.51%  "
$ws.Range("D29").Value = "'5.216"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.12%  "
$ws.Range("D30").Value = "'116.08"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.67%  "
$ws.Range("E31").Value = "  -1.84%  "
$ws.Range("D32").Value = "'1.193"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Value = "'0.7533"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.99%  "
$ws.Range("D34").Value = "'4.426"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.60%  "
$ws.Range("D35").Value = "'2.804"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.45%  "
$ws.Range("D36").Value = "'1.005"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.21%  "
$ws.Range("D37").Value = "'1.110"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.18%  "
$ws.Range("E38").Value = "  -0.47%  "
$ws.Range("D39").Value = "'0.05253"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.32%  "
$ws.Range("D40").Value = "'7.144"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.00%  "
$ws.Range("D41").Value = "'2.862"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.78%  "
$ws.Range("D42").Value = "'0.1687"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.24%  "
$ws.Range("D43").Value = "'0.5033"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.51%  "
$ws.Range("D44").Value = "'8.607"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.10%  "
$ws.Range("E45").Value = "  -1.38%  "
$ws.Range("D46").Value = "'106.97"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.51%  "
$ws.Range("D47").Value = "'0.4737"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.59%  "
$ws.Range("D48").Value = "'1.006"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.30%  "
$ws.Range("E49").Value = "  -1.84%  "
$ws.Range("D50").Value = "'1.654"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.69%  "
$ws.Range("D51").Value = "'1.800"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.84%  "
